$d = $word.ActiveDocument

# 1. Replace the UC3_Service_Date_Report text with UC7_View_Available_Coupons (split into 3 runs)
$found = $d.Content.Find.Execute("UC3_Service_Date_Report  (Need to rename RTM)", $true, $false, $false, $false, $false, $true, 1, $false, "UC7_View_Available_Coupons", 2)
